# Update timing statistics (Tiempo_Mínimo, Tiempo_Máximo, Tiempo_Promedio)
# for rows 2-4 on the "Data" sheet, per new experiment run results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2
$ws.Range("E2").Value = 0.00018558
$ws.Range("F2").Value = 0.01528371
$ws.Range("G2").Value = 0.00035320563900000004

# Row 3
$ws.Range("E3").Value = 0.00192375
$ws.Range("F3").Value = 0.00979911
$ws.Range("G3").Value = 0.00244350859030837

# Row 4
$ws.Range("E4").Value = 0.00434079
$ws.Range("F4").Value = 0.01065708
$ws.Range("G4").Value = 0.005334742019230769
